# "Changed locator and worked on author"
#
# tc002!A2 held the old "locator" string ("STG- PulseCodeOnAzureCloud");
# the author updated it to a new project/locator name, which lands as a
# brand-new shared-string entry ("STG- SPARK Modernization"), and the
# active/selected sheet+cell moved from tc010!G5 to tc002!A7 (tabSelected
# hops from tc010 to tc002 and the active cell on tc002 becomes A7).

$wb = $excel.ActiveWorkbook

# Update the locator text on tc002.
$ws = $wb.Worksheets.Item("tc002")
$ws.Range("A2").Value = "STG- SPARK Modernization"

# Make tc002 the active sheet / tab, with A7 selected - mirrors the
# workbook-level active tab + sheetView selection move in the diff
# (tc010 loses tabSelected, tc002 gains it and its selection becomes A7).
$ws.Activate()
$ws.Range("A7").Select()
